$d = $word.ActiveDocument

# --- Paragraph 3 ("Here is an example of their usage...") ---
# Remove run-level yellow highlighting.
$p3 = $d.Paragraphs.Item(3)
$p3.Range.HighlightColorIndex = 0

# The paragraph-mark's own run-properties (pPr/rPr) still carries the
# yellow highlight and isn't reachable through Range.HighlightColorIndex,
# so reset it by deleting + reinserting the mark (the new mark inherits
# the current, non-highlighted formatting).
$markPos = $p3.Range.End - 1
$markRange = $d.Range($markPos, $p3.Range.End)
$markRange.Delete()
$d.Range($markPos, $markPos).InsertParagraphAfter()

# --- Paragraph 4 ("From the first pass, it is now possible...") ---
# Replace its content: drop the highlight, and rewrite the text after
# "given that a" with the expanded Robustness/Sequence-diagram wording,
# relocating the _GoBack bookmark to the end of the paragraph.
$p4 = $d.Paragraphs.Item(4)
$rng = $d.Range($p4.Range.Start, $p4.Range.End - 1)
$xmlFrag = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="1E51B09D" w14:textId="390FE3D5" w:rsidR="00102AC5" w:rsidRDefault="002B4863" w:rsidP="00102AC5"><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">From the first pass, it is now possible to determine </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">the </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">properties/methods of a class, </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>given that a</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> general overview of their responsibilities has been provided. </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Note that each class has an ID (which is used as a primary-key, in the project’s database). Considering the Member class, The Game Café wants to know the </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>details for a Member (e.g. Membership Type)</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>A Member is linked to the Booking class (as they can make multiple bookings) and the eSports-Event class (as they can also purchase tickets for an eSports Event).</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>The Booking class will store all of the relevant details for a Booking (such as which Member has made that Booking)</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Hence, they are linked to the Member class, as well as the Hardware class (as Member’s can choose a </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>particular platform</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> for their Booking)</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Considering the Hardware class, there is </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">a </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Software</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> class, for the Software that runs on </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>that piece of Hardware (if appropriate)</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>. This class stores the appropriate details for that piece of Software.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng.InsertXML($xmlFrag)

Write-Output "edit applied"
